# Strip the document-defaults overrides back down to the minimal set,
# matching the "download tc, tcn, and tl files from GD" commit:
#   rPrDefault keeps only rFonts / sz / szCs / lang
#   pPrDefault keeps only spacing (line/lineRule)
#
# Word's COM model has no direct "docDefaults" object, so we round-trip
# the whole package through Document.WordOpenXML, patch the <w:docDefaults>
# block with a regex (tolerant of attribute-order / boolean-shorthand
# differences in the serializer), and assign the patched XML back.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

$pattern = '<w:docDefaults>.*?</w:docDefaults>'
$replacement = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="de"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

$newXml = [System.Text.RegularExpressions.Regex]::Replace($xml, $pattern, $replacement, [System.Text.RegularExpressions.RegexOptions]::Singleline)

$d.WordOpenXML = $newXml
